$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value to a cell while preserving its original
# (default) cell style, so only the value itself changes.
function Set-TextValue($cell, $value) {
    $range = $ws.Range($cell)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue "D2" "261.50"
Set-TextValue "D3" "26.38"
Set-TextValue "E3" "-3.74%"
Set-TextValue "D4" "4.707"
Set-TextValue "E4" "0.55%"
Set-TextValue "D5" "0.06162"
Set-TextValue "E5" "1.13%"
Set-TextValue "D6" "6.709"
Set-TextValue "E6" "0.66%"
Set-TextValue "D7" "0.8509"
Set-TextValue "E7" "0.24%"
Set-TextValue "D8" "0.9113"
Set-TextValue "E8" "-1.40%"
Set-TextValue "E9" "0.08%"
Set-TextValue "D10" "0.05202"
Set-TextValue "E10" "6.93%"
Set-TextValue "D11" "0.07097"
Set-TextValue "E11" "-0.01%"
Set-TextValue "D12" "0.03120"
Set-TextValue "E12" "1.41%"
Set-TextValue "D13" "0.09046"
Set-TextValue "E13" "-0.16%"
Set-TextValue "D14" "0.001535"
Set-TextValue "E14" "-0.07%"
Set-TextValue "D15" "0.0006182"
Set-TextValue "E15" "1.43%"
Set-TextValue "D16" "0.005958"
Set-TextValue "E16" "-2.57%"
Set-TextValue "E17" "0.04%"
Set-TextValue "D18" "3.173"
Set-TextValue "E18" "0.78%"
Set-TextValue "E19" "1.12%"
Set-TextValue "D21" "0.1300"
Set-TextValue "E21" "-0.78%"
Set-TextValue "D22" "4.090"
Set-TextValue "E22" "0.08%"
Set-TextValue "D23" "0.04245"
Set-TextValue "E23" "0.16%"
Set-TextValue "D24" "0.001177"
Set-TextValue "E25" "6.50%"
Set-TextValue "E26" "0.04%"
Set-TextValue "D40" "0.03998"
Set-TextValue "E40" "3.61%"
Set-TextValue "D41" "0.1112"
Set-TextValue "E41" "0.04%"
Set-TextValue "D42" "0.004128"
Set-TextValue "E42" "1.41%"
Set-TextValue "E43" "-3.34%"
Set-TextValue "D44" "0.01328"
Set-TextValue "E44" "-18.16%"
Set-TextValue "E45" "0.40%"
Set-TextValue "E46" "0.06%"
Set-TextValue "E48" "87.55%"
Set-TextValue "E49" "0.06%"
Set-TextValue "D50" "0.0002002"
Set-TextValue "E50" "0.06%"
